# Insert a new weekly price record for Coliflor at row 384, pushing the
# existing rows 384-431 down to 385-432 (dimension grows from R431 to R432).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 384..431 down by one row, leaving row 384 empty (but carrying
# the date-formatted style from the row above onto column D, same as Excel
# does natively on a row insert).
$ws.Rows.Item(384).Insert()

# Populate the newly inserted row 384 with the new record.
$ws.Range("A384").Value = 5
$ws.Range("B384").Value = "Macroferia Regional de Talca"
$ws.Range("C384").Value = "Maule"
$ws.Range("D384").Value = 45077
$ws.Range("E384").Value = 7
$ws.Range("F384").Value = 100112008
$ws.Range("G384").Value = "Coliflor"
$ws.Range("H384").Value = "Sin especificar"
$ws.Range("I384").Value = "Primera"
$ws.Range("J384").Value = 5000
$ws.Range("K384").Value = 700
$ws.Range("L384").Value = 700
$ws.Range("M384").Value = 700
$ws.Range("N384").Value = "`$/unidad"
$ws.Range("O384").Value = "Región del Maule"
$ws.Range("P384").Value = 700
$ws.Range("Q384").Value = 1
$ws.Range("R384").Value = "Hortaliza"
